$wb = $excel.ActiveWorkbook

# Hunk @ diff line 2628 -- sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3500.2727
$ws.Cells.Item(40, 9).Value = 6250.5
$ws.Cells.Item(40, 10).Value = 2889.111
$ws.Cells.Item(40, 11).Value = 6250.5
$ws.Cells.Item(40, 12).Value = 2889.111
$ws.Cells.Item(40, 13).Value = -6075.5
$ws.Cells.Item(40, 14).Value = -3239.111

# Hunk @ diff line 3228 -- sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(52, 8).Value = 1743.3334
$ws.Cells.Item(52, 10).Value = 3000
$ws.Cells.Item(52, 12).Value = 9000
$ws.Cells.Item(52, 14).Value = -9320

# Hunk @ diff line 4082 -- sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 4624.75
$ws.Cells.Item(69, 9).Value = 4998.25
$ws.Cells.Item(69, 10).Value = 4251.25
$ws.Cells.Item(69, 11).Value = 14994.75
$ws.Cells.Item(69, 12).Value = 12753.75
$ws.Cells.Item(69, 13).Value = -14120.75
$ws.Cells.Item(69, 14).Value = -14501.75

# Hunk @ diff line 4235 -- sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 4624.75
$ws.Cells.Item(72, 9).Value = 4998.25
$ws.Cells.Item(72, 10).Value = 4251.25
$ws.Cells.Item(72, 11).Value = 44984.25
$ws.Cells.Item(72, 12).Value = 38261.25
$ws.Cells.Item(72, 13).Value = -40616.25
$ws.Cells.Item(72, 14).Value = -46997.25

# Hunk @ diff line 7256 -- sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 7403.7554
$ws.Cells.Item(132, 9).Value = 6838.967
$ws.Cells.Item(132, 10).Value = 8533.333000000001
$ws.Cells.Item(132, 11).Value = 20516.901
$ws.Cells.Item(132, 12).Value = 25599.999
$ws.Cells.Item(132, 13).Value = -17986.901
$ws.Cells.Item(132, 14).Value = -30659.999

# Hunk @ diff line 7501 -- sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1441.9667
$ws.Cells.Item(137, 9).Value = 1403.8372
$ws.Cells.Item(137, 10).Value = 1538.4117
$ws.Cells.Item(137, 11).Value = 4211.5116
$ws.Cells.Item(137, 12).Value = 4615.2351
$ws.Cells.Item(137, 13).Value = -1661.5116
$ws.Cells.Item(137, 14).Value = -9715.2351

# Hunk @ diff line 11365 -- sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1641.405
$ws.Cells.Item(74, 9).Value = 1612.7102
$ws.Cells.Item(74, 10).Value = 1839.4
$ws.Cells.Item(74, 11).Value = 1612.7102
$ws.Cells.Item(74, 12).Value = 1839.4
$ws.Cells.Item(74, 13).Value = -738.7102
$ws.Cells.Item(74, 14).Value = -3587.4

# Hunk @ diff line 11512 -- sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1641.405
$ws.Cells.Item(77, 9).Value = 1612.7102
$ws.Cells.Item(77, 10).Value = 1839.4
$ws.Cells.Item(77, 11).Value = 8063.550999999999
$ws.Cells.Item(77, 12).Value = 9197
$ws.Cells.Item(77, 13).Value = -3695.550999999999
$ws.Cells.Item(77, 14).Value = -17933

# Hunk @ diff line 12734 -- sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2010
$ws.Cells.Item(102, 9).Value = 2010
$ws.Cells.Item(102, 11).Value = 2010
$ws.Cells.Item(102, 13).Value = -388

# Hunk @ diff line 13705 -- sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1852.7778
$ws.Cells.Item(122, 9).Value = 1775.7333
$ws.Cells.Item(122, 10).Value = 2238
$ws.Cells.Item(122, 11).Value = 5327.199900000001
$ws.Cells.Item(122, 12).Value = 6714
$ws.Cells.Item(122, 13).Value = -2877.199900000001
$ws.Cells.Item(122, 14).Value = -11614

# Hunk @ diff line 14189 -- sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3327.7812
$ws.Cells.Item(132, 9).Value = 1415.3778
$ws.Cells.Item(132, 10).Value = 7857.1577
$ws.Cells.Item(132, 11).Value = 4246.1334
$ws.Cells.Item(132, 12).Value = 23571.4731
$ws.Cells.Item(132, 13).Value = -1716.1334
$ws.Cells.Item(132, 14).Value = -28631.4731

# Hunk @ diff line 15428 -- sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()

# Hunk @ diff line 15670 -- sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1573.4642
$ws.Cells.Item(20, 9).Value = 579.4666999999999
$ws.Cells.Item(20, 10).Value = 2720.3845
$ws.Cells.Item(20, 11).Value = 579.4666999999999
$ws.Cells.Item(20, 12).Value = 2720.3845
$ws.Cells.Item(20, 13).Value = -332.4666999999999
$ws.Cells.Item(20, 14).Value = -3214.3845

# Hunk @ diff line 16393 -- sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 47000
$ws.Cells.Item(35, 10).Value = 49666.668
$ws.Cells.Item(35, 12).Value = 49666.668
$ws.Cells.Item(35, 14).Value = -50286.668

# Hunk @ diff line 19517 -- sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2093.8333
$ws.Cells.Item(99, 9).Value = 1563.5
$ws.Cells.Item(99, 10).Value = 3950
$ws.Cells.Item(99, 11).Value = 1563.5
$ws.Cells.Item(99, 12).Value = 3950
$ws.Cells.Item(99, 13).Value = -65.5
$ws.Cells.Item(99, 14).Value = -6946

# Hunk @ diff line 19817 -- sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2647.7778
$ws.Cells.Item(105, 9).Value = 2150
$ws.Cells.Item(105, 10).Value = 3270
$ws.Cells.Item(105, 11).Value = 2150
$ws.Cells.Item(105, 12).Value = 3270
$ws.Cells.Item(105, 13).Value = -403
$ws.Cells.Item(105, 14).Value = -6764

# Hunk @ diff line 21220 -- sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4488.385
$ws.Cells.Item(134, 9).Value = 2238.2856
$ws.Cells.Item(134, 10).Value = 6012.645
$ws.Cells.Item(134, 11).Value = 6714.8568
$ws.Cells.Item(134, 12).Value = 18037.935
$ws.Cells.Item(134, 13).Value = -4179.8568
$ws.Cells.Item(134, 14).Value = -23107.935

# Hunk @ diff line 23154 -- sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4017484.8
$ws.Cells.Item(31, 9).Value = 1033.7031
$ws.Cells.Item(31, 10).Value = 17546582
$ws.Cells.Item(31, 11).Value = 1033.7031
$ws.Cells.Item(31, 12).Value = 17546582
$ws.Cells.Item(31, 13).Value = -738.7030999999999
$ws.Cells.Item(31, 14).Value = -17547172

# Hunk @ diff line 23307 -- sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 4017484.8
$ws.Cells.Item(34, 9).Value = 1033.7031
$ws.Cells.Item(34, 10).Value = 17546582
$ws.Cells.Item(34, 11).Value = 1033.7031
$ws.Cells.Item(34, 12).Value = 17546582
$ws.Cells.Item(34, 13).Value = -831.7030999999999
$ws.Cells.Item(34, 14).Value = -17546986

# Hunk @ diff line 31397 -- sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(54, 8).Value = 2839.3845
$ws.Cells.Item(54, 10).Value = 3000
$ws.Cells.Item(54, 12).Value = 9000
$ws.Cells.Item(54, 14).Value = -10118

# Hunk @ diff line 32774 -- sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 2839.6155
$ws.Cells.Item(81, 10).Value = 2839.6155
$ws.Cells.Item(81, 12).Value = 8518.8465
$ws.Cells.Item(81, 14).Value = -10764.8465

# Hunk @ diff line 32927 -- sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(84, 8).Value = 2839.6155
$ws.Cells.Item(84, 10).Value = 2839.6155
$ws.Cells.Item(84, 12).Value = 25556.5395
$ws.Cells.Item(84, 14).Value = -36788.5395

# Hunk @ diff line 39502 -- sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(74, 8).Value = 19966.666
$ws.Cells.Item(74, 10).Value = 19966.666
$ws.Cells.Item(74, 12).Value = 19966.666
$ws.Cells.Item(74, 14).Value = -21838.666

# Hunk @ diff line 39646 -- sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(77, 8).Value = 19966.666
$ws.Cells.Item(77, 10).Value = 19966.666
$ws.Cells.Item(77, 12).Value = 59899.99800000001
$ws.Cells.Item(77, 14).Value = -69259.99800000001

# Hunk @ diff line 41888 -- sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 10348.223
$ws.Cells.Item(123, 10).Value = 10348.223
$ws.Cells.Item(123, 12).Value = 10348.223
$ws.Cells.Item(123, 14).Value = -15248.223

# Hunk @ diff line 42332 -- sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1192658.5
$ws.Cells.Item(132, 9).Value = 3474119.8
$ws.Cells.Item(132, 10).Value = 2330.9565
$ws.Cells.Item(132, 11).Value = 10422359.4
$ws.Cells.Item(132, 12).Value = 6992.869499999999
$ws.Cells.Item(132, 13).Value = -10419829.4
$ws.Cells.Item(132, 14).Value = -12052.8695

# Hunk @ diff line 55768 -- sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3211.75
$ws.Cells.Item(122, 9).Value = 2724.9092
$ws.Cells.Item(122, 10).Value = 4672.273
$ws.Cells.Item(122, 11).Value = 8174.7276
$ws.Cells.Item(122, 12).Value = 14016.819
$ws.Cells.Item(122, 13).Value = -5724.7276
$ws.Cells.Item(122, 14).Value = -18916.819
